$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1099.7142
$ws.Range("I62").Value = 1019.6
$ws.Range("K62").Value = 1019.6
$ws.Range("M62").Value = -395.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1099.7142
$ws.Range("I65").Value = 1019.6
$ws.Range("K65").Value = 5098
$ws.Range("M65").Value = -1978

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1535.5555
$ws.Range("I70").Value = 1240
$ws.Range("J70").Value = 1572.5
$ws.Range("K70").Value = 3720
$ws.Range("L70").Value = 4717.5
$ws.Range("M70").Value = -3450
$ws.Range("N70").Value = -5257.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1535.5555
$ws.Range("I73").Value = 1240
$ws.Range("J73").Value = 1572.5
$ws.Range("K73").Value = 3720
$ws.Range("L73").Value = 4717.5
$ws.Range("M73").Value = -2784
$ws.Range("N73").Value = -6589.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5984.615
$ws.Range("I76").Value = 3560
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 3560
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = -3245
$ws.Range("N76").Value = -8130

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5984.615
$ws.Range("I79").Value = 3560
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 3560
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = -2468
$ws.Range("N79").Value = -9684

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6667538
$ws.Range("I137").Value = 827.15
$ws.Range("J137").Value = 20000960
$ws.Range("K137").Value = 2481.45
$ws.Range("L137").Value = 60002880
$ws.Range("M137").Value = 68.55000000000018
$ws.Range("N137").Value = -60007980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 78.75
$ws.Range("I5").Value = 46
$ws.Range("J5").Value = 133.33333
$ws.Range("K5").Value = 46
$ws.Range("L5").Value = 133.33333
$ws.Range("M5").Value = 66
$ws.Range("N5").Value = -357.33333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1325.64
$ws.Range("I110").Value = 1116.8667
$ws.Range("J110").Value = 1638.8
$ws.Range("K110").Value = 1116.8667
$ws.Range("L110").Value = 1638.8
$ws.Range("M110").Value = 928.1333
$ws.Range("N110").Value = -5728.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6946479.5
$ws.Range("I132").Value = 8066224
$ws.Range("J132").Value = 4062.8
$ws.Range("K132").Value = 24198672
$ws.Range("L132").Value = 12188.4
$ws.Range("M132").Value = -24196142
$ws.Range("N132").Value = -17248.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 78.75
$ws.Range("I4").Value = 46
$ws.Range("J4").Value = 133.33333
$ws.Range("K4").Value = 46
$ws.Range("L4").Value = 133.33333
$ws.Range("M4").Value = 69
$ws.Range("N4").Value = -363.33333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2422.3044
$ws.Range("I134").Value = 1608.6061
$ws.Range("J134").Value = 4487.846
$ws.Range("K134").Value = 4825.8183
$ws.Range("L134").Value = 13463.538
$ws.Range("M134").Value = -2290.8183
$ws.Range("N134").Value = -18533.538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1058.9
$ws.Range("I16").Value = 1020.125
$ws.Range("J16").Value = 1214
$ws.Range("K16").Value = 1020.125
$ws.Range("L16").Value = 1214
$ws.Range("M16").Value = -733.125
$ws.Range("N16").Value = -1788

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6948579
$ws.Range("I31").Value = 4255.282
$ws.Range("K31").Value = 4255.282
$ws.Range("M31").Value = -3960.282

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6948579
$ws.Range("I34").Value = 4255.282
$ws.Range("K34").Value = 4255.282
$ws.Range("M34").Value = -4053.282

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1058.9
$ws.Range("I113").Value = 1020.125
$ws.Range("J113").Value = 1214
$ws.Range("K113").Value = 1020.125
$ws.Range("L113").Value = 1214
$ws.Range("M113").Value = 1149.875
$ws.Range("N113").Value = -5554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 13160201
$ws.Range("I132").Value = 20002054
$ws.Range("J132").Value = 2789.5386
$ws.Range("K132").Value = 60006162
$ws.Range("L132").Value = 8368.6158
$ws.Range("M132").Value = -60003632
$ws.Range("N132").Value = -13428.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 568.3333
$ws.Range("I122").Value = 570.0714
$ws.Range("K122").Value = 5130.6426
$ws.Range("M122").Value = -2680.6426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 880.8461
$ws.Range("I132").Value = 707.06665
$ws.Range("J132").Value = 1117.8182
$ws.Range("K132").Value = 6363.59985
$ws.Range("L132").Value = 10060.3638
$ws.Range("M132").Value = -3833.59985
$ws.Range("N132").Value = -15120.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2779545
$ws.Range("I122").Value = 3922717.2
$ws.Range("J122").Value = 3270
$ws.Range("K122").Value = 11768151.6
$ws.Range("L122").Value = 9810
$ws.Range("M122").Value = -11765701.6
$ws.Range("N122").Value = -14710

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 52780
$ws.Range("J128").Value = 52780
$ws.Range("L128").Value = 52780
$ws.Range("N128").Value = -62740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2827.2097
$ws.Range("I132").Value = 2129.3408
$ws.Range("J132").Value = 4533.1113
$ws.Range("K132").Value = 6388.0224
$ws.Range("L132").Value = 13599.3339
$ws.Range("M132").Value = -3858.0224
$ws.Range("N132").Value = -18659.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 962.96155
$ws.Range("I22").Value = 484
$ws.Range("J22").Value = 1441.9231
$ws.Range("K22").Value = 484
$ws.Range("L22").Value = 1441.9231
$ws.Range("M22").Value = -189
$ws.Range("N22").Value = -2031.9231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 962.96155
$ws.Range("I27").Value = 484
$ws.Range("J27").Value = 1441.9231
$ws.Range("K27").Value = 484
$ws.Range("L27").Value = 1441.9231
$ws.Range("M27").Value = -377
$ws.Range("N27").Value = -1655.9231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9165
$ws.Range("I40").Value = 10025
$ws.Range("K40").Value = 10025
$ws.Range("M40").Value = -9889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5987.7144
$ws.Range("J122").Value = 4275.5557
$ws.Range("L122").Value = 12826.6671
$ws.Range("N122").Value = -17726.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 87678
$ws.Range("J131").Value = 87678
$ws.Range("L131").Value = 87678
$ws.Range("N131").Value = -97758
